$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the static, per-row DOI citations in AP10:AP39 with formulas
#     that reference the previous row (AP10 -> =AP9, AP11 -> =AP10, etc.),
#     matching the "fill down from AP9" edit made by the author. Assigning
#     the whole block in one shot lets the engine coalesce AP11:AP39 into a
#     single shared-formula group, exactly like Excel's own fill-down does.
$ws.Range("AP10").Formula = "=AP9"
$ws.Range("AP11:AP39").Formula = "=AP10"

# --- The old per-row hyperlink that covered AP10:AP39 no longer matches the
#     (now-formula-driven) content, so drop it; the AP9 hyperlink is untouched.
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $hl = $ws.Hyperlinks.Item($i)
    if ($hl.Range.Row -eq 10) {
        $hl.Delete()
    }
}

# --- Widen the citation-note columns (AL, AM) that were touched in the same
#     editing pass.
$ws.Columns.Item(38).ColumnWidth = 18.833333333333332
$ws.Columns.Item(39).ColumnWidth = 29.5

# --- Update the visible selection to the newly-edited range.
$ws.Range("AP10:AP39").Select()
